# Update the "想去人数" (want-to-go count) figures that changed between
# the previous gh-pages data snapshot and the new one generated at 456a3b4.
#
# The same underlying rows are duplicated on the "展览" sheet and on the
# aggregated "全部类型" sheet, so both need to be updated identically.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 715
    7  = 23
    11 = 4589
    12 = 4415
    15 = 151
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
